# Trade #76 closed at 2026-02-17 21:13:24 - unknown UNKNOWN +0.000%
#
# This script:
#  1. Updates the "Summary" sheet aggregate stats (Total Trades, Win Rate %)
#  2. Updates the "Strategy Status" sheet row for MarketMaking (Trades, Win Rate %)
#  3. Closes the open MarketMaking trade (row 105 in "All Trades", row 72 in
#     "MarketMaking") that was still OPEN, marking it CLOSED with an exit price,
#     exit reason and duration.
#  4. Appends a brand-new OPEN trade row to both "All Trades" and "MarketMaking".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 104
$summary.Range("B9").Value = 47.12

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet (row 5 = MarketMaking)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value = 71
$status.Range("G5").Value = 49.3

# ---------------------------------------------------------------------------
# 3) All Trades sheet - close trade #104 (row 105)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G105").Value = 0.03
$allTrades.Range("H105").Value = "CLOSED"
$allTrades.Range("K105").Value = 101.15
$allTrades.Range("L105").Value = "early_exit"
$allTrades.Range("M105").Value = 0.11

# Append new open trade (#137) as row 138
$allTrades.Range("A138").Value = 137
$allTrades.Range("B138").NumberFormat = "@"
$allTrades.Range("B138").Value = "2026-02-17"
$allTrades.Range("C138").Value = "21:13:18"
$allTrades.Range("D138").Value = "MarketMaking"
$allTrades.Range("E138").Value = "DOWN"
$allTrades.Range("F138").Value = 0.03
$allTrades.Range("H138").Value = "OPEN"
$allTrades.Range("I138").Value = 0
$allTrades.Range("J138").Value = 0
$allTrades.Range("K138").Value = 101.1496151053151
$allTrades.Range("M138").Value = 0
$allTrades.Range("N138").Value = 0
$allTrades.Range("O138").Value = 0
$allTrades.Range("P138").Value = 0.6
$allTrades.Range("Q138").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# 4) MarketMaking sheet - close trade #104 (row 72)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G72").Value = 0.03
$marketMaking.Range("H72").Value = "CLOSED"
$marketMaking.Range("K72").Value = 101.15
$marketMaking.Range("P72").Value = "early_exit"
$marketMaking.Range("Q72").Value = 0.11

# Append new open trade (#137) as row 105
$marketMaking.Range("A105").Value = 137
$marketMaking.Range("B105").NumberFormat = "@"
$marketMaking.Range("B105").Value = "2026-02-17"
$marketMaking.Range("C105").Value = "21:13:18"
$marketMaking.Range("D105").Value = "MarketMaking"
$marketMaking.Range("E105").Value = "DOWN"
$marketMaking.Range("F105").Value = 0.03
$marketMaking.Range("H105").Value = "OPEN"
$marketMaking.Range("I105").Value = 0
$marketMaking.Range("J105").Value = 0
$marketMaking.Range("K105").Value = 101.1496151053151
$marketMaking.Range("L105").Value = 0
$marketMaking.Range("M105").Value = 0
$marketMaking.Range("N105").Value = 0.6
$marketMaking.Range("O105").Value = "Normal spread capture: 19600 bps"
$marketMaking.Range("Q105").Value = 0
